$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Search")

# Insert a new row at position 4, shifting the existing row 4 (and below) down
$ws.Rows.Item(4).Insert()

# Populate the newly inserted row 4 with the search_full_text labels
$ws.Range("A4").Value = '${msg.getProperty(''search_full_text'')}'
$ws.Range("B4").Value = '${search_full_text}'
